$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new comment rows appended to the table (rows 19 and 20), using the
# next available shared-string entries.
$ws.Range("A19").Value = "*SortedDictionary does not support Range operation."
$ws.Range("A20").Value = "*SortedSet is faster than OrderedSet"

# New rows get the existing "note" look, but in red text (new font/cellXf).
$ws.Range("A19:A20").Font.Color = 255

# Reflect the current on-screen selection at save time.
[void]$ws.Range("C23").Select()

# Print setup was touched (paper size set explicitly to Letter/A4 = 9).
$ws.PageSetup.PaperSize = 9
